$wb = $excel.ActiveWorkbook

# Sheet "展览" (展览信息) - first data sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 632
$ws1.Range("F5").Value = 177
$ws1.Range("F6").Value = 9507
$ws1.Range("F10").Value = 1182
$ws1.Range("F11").Value = 153
$ws1.Range("F12").Value = 103
$ws1.Range("F13").Value = 21
$ws1.Range("F17").Value = 257
$ws1.Range("F18").Value = 1305

# Sheet "全部类型" (combined list) - rows offset by 1 vs "展览"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 632
$ws4.Range("F6").Value = 177
$ws4.Range("F7").Value = 9507
$ws4.Range("F11").Value = 1182
$ws4.Range("F12").Value = 153
$ws4.Range("F13").Value = 103
$ws4.Range("F14").Value = 21
$ws4.Range("F18").Value = 257
$ws4.Range("F19").Value = 1305
